$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

for ($r = 2; $r -le 157; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.EndsWith(", System")) {
        $parts = $val -split ", "
        $last = $parts[-1]
        $rest = $parts[0..($parts.Length - 2)]
        $newVal = @($last) + $rest -join ", "
        $cell.Value2 = $newVal
    }
}
